# "add minor changes to excel"
#
# 1) Reorder worksheets: move "DateToText" so it comes before "LeftRight"
#    (previously LeftRight, then DateToText).
# 2) Make "Substitute" the active/selected sheet (tab).
# 3) On "DateToText": add a "day"/"year" breakdown using ANCHORARRAY off the
#    existing J (text date) dynamic array column, replacing the old helper
#    label in K and the intermediate re-TEXT()'d column in L, and dropping
#    the now-unused M column.

$wb = $excel.ActiveWorkbook

# --- 1. Reorder sheets: DateToText before LeftRight -----------------------
$dateToText = $wb.Worksheets.Item("DateToText")
$leftRight  = $wb.Worksheets.Item("LeftRight")
$dateToText.Move($leftRight)

# --- 2. Rework the day/year helper columns on DateToText ------------------
$ws = $wb.Worksheets.Item("DateToText")

# Clear the old K (label) and M (RIGHT-of-L year) columns entirely.
$ws.Range("K1:K10").Clear()
$ws.Range("M1:M10").Clear()

# New headers
$ws.Range("K1").Value = "day"
$ws.Range("L1").Value = "year"

# K2:K10 <- day part of the J (TEXT) dynamic array; L2:L10 <- year part.
$ws.Range("K2:K10").FormulaArray = "=LEFT(J2:J10,2)"
$ws.Range("L2:L10").FormulaArray = "=RIGHT(J2:J10,4)"

# Selection bookkeeping on the edited sheet (done while it's still active).
$ws.Range("L17").Select()

# --- 3. Active tab becomes "Substitute" ------------------------------------
$substitute = $wb.Worksheets.Item("Substitute")
$substitute.Activate()
